$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 'Lea este sitio en español: <a href=''https://smartenies.shinyapps.io/ges_health_study_app_espanol''> Mapas interactivos del Estudio de salud comunitaria del GES</a>'
$ws.Range("C2").Value = 'Read this site in English: <a href=''https://smartenies.shinyapps.io/ges_health_study_app_english''> GES Community Health Study Interactive Maps</a>'
$ws.Range("B4").Value = 'GES Community Health Study Interactive Maps'
$ws.Range("C4").Value = 'Mapas interactivos del Estudio de salud comunitaria del GES'
$ws.Range("B7").Value = 'Welcome to the GES Community Health Study Interactive Maps app! We are excited to share this important tool with you. The Interactive Maps app was developed to help community members learn more about the environmental conditions that currently and historically exist in the Globeville, Elyria, and Swansea neighborhoods of Denver.'
$ws.Range("C7").Value = '¡Bienvenido a la aplicación de mapas interactivos del Estudio de salud comunitaria del GES! Nos complace compartir esta importante herramienta con usted. La aplicación de mapas interactivos se desarrolló para ayudar a los miembros de la comunidad a aprender más sobre las condiciones ambientales que existen actualmente e históricamente en los vecindarios de Globeville, Elyria y Swansea de Denver.'
$ws.Range("B10").Value = 'We invite you to continue reading to learn more about the GES Community Health Study Interactive Maps app.'
$ws.Range("B45").Value = 'More information about this map and our findings can be found on the <a href=''https://www.geshealthstudy.org/maps/community-mapping-locations''> GES Community Health Study website</a>.'
$ws.Range("B48").Value = 'In addition to our community mapping exercise, the GES Community Health Study has generated a number of maps documenting environmental and health conditions in the GES neighborhoods. We have included these maps here for you to explore.'
$ws.Range("C48").Value = 'Además de nuestro ejercicio de mapeo comunitario, el Estudio de salud comunitaria del GES ha generado una serie de mapas que documentan las condiciones ambientales y de salud en los vecindarios de GES. Hemos incluido estos mapas aquí para que los explore.'
$ws.Range("C70").Value = 'Esta herramienta de mapeo dinámico es una colaboración entre ENVIRONS y el Consejo Comunitario del Estudio de salud comunitaria del GES. Esta herramienta permite a los usuarios comparar las características de la población y las exposiciones ambientales en los vecindarios de Denver. Se creó teniendo en cuenta las necesidades de datos de la comunidad y se basó en los comentarios de nuestro Consejo Comunitario de GES.'
$ws.Range("C74").Value = 'Puede encontrar más información sobre el Consejo Comunitario de GES, el equipo ENVRONS y el Comité Directivo del Estudio de salud comunitaria del GES en nuestro sitio web: <a href=''https://www.geshealthstudy.org/about/who-we-are''>Estudio de Salud Comunitaria de GES</a>.'

$ws.Range("B7").Select()
